$d = $word.ActiveDocument

# --- 1. Replace "EXT. FINISH LINE - DAY" with 4 runs forming "INT. GE BOSTON OFFICE - DAY" ---
$p1 = $d.Paragraphs.Item(1)
$start1 = $p1.Range.Start
$end1 = $p1.Range.End

$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>INT</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>GE BOSTON OFFICE</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> - DAY</w:t></w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$r1 = $d.Range($start1, $end1)
$r1.InsertXML($xml1)

# --- 2. Replace "HARRIET" with "NARRATOR" ---
$d.Content.Find.Execute("HARRIET", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NARRATOR", 2)

# --- 3. Replace the dialogue block (paragraphs: "I can’t believe...", blank, "PATRICIA",
#         "Slow and steady...") and the trailing bookmark-only paragraph with the new narration
#         paragraph, preserving the _GoBack bookmark in the middle of the new runs.
$p4 = $d.Paragraphs.Item(4)
$start3 = $p4.Range.Start
$end3 = $d.Content.End

$xml3 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
<w:r><w:t>SCHMENDAN</w:t></w:r>
<w:r><w:t xml:space="preserve"> was unfortunately, but rightfully, fired from GE and went on to regret his shameful display for the rest of his life.</w:t></w:r>
<w:r><w:t xml:space="preserve"> On the other hand, </w:t></w:r>
<w:r><w:t xml:space="preserve">SCHMORDAN went on to get her </w:t></w:r>
<w:r><w:t>dream job in GE’s Boston office and was highly regarded for her</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t xml:space="preserve"> poise for the entirety of her career. </w:t></w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$r3 = $d.Range($start3, $end3)
$r3.InsertXML($xml3)

Write-Output "edit complete"
